$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Core content change: A7 "threshold" -> "binarize"
$ws.Range("A7").Value = "binarize"

# Column width changes
$ws.Columns.Item(1).ColumnWidth = 14.140625
$ws.Columns.Item(2).ColumnWidth = 20.7109375

# Selection change
$ws.Range("E12:E13").Select()
